$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Tasks completed this week" / "Tasks to complete next week" entries.
# A19 previously shared the same string as B19 ("High Fidelity Prototype").
# Now A19 gets new text, B19 gets different new text, and A20 gets a new line too.
$ws.Range("A19").Value = "Started building frotend using Processing."
$ws.Range("A20").Value = "Wrote basics for our own library for UI."
$ws.Range("B19").Value = "Further develop the frontend."

# Re-select a different active cell, matching the new selection in the sheet view.
$ws.Range("C19").Select()
